$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1531
$ws.Cells.Item(1531,1).Value = 1530
$ws.Cells.Item(1531,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1531,3).Value = "3:55 PM"
$ws.Cells.Item(1531,4).Value = "LO521"
$ws.Cells.Item(1531,5).Value = "Prague"
$ws.Cells.Item(1531,6).Value = "(PRG)"
$ws.Cells.Item(1531,7).Value = "LOT "
$ws.Cells.Item(1531,8).Value = "E75S"
$ws.Cells.Item(1531,9).Value = "(SP-LIL)"
$ws.Cells.Item(1531,10).Value = "4:10 PM"
$ws.Cells.Item(1531,12).Value = "0 hours, 15 minutes"

# Row 1532
$ws.Cells.Item(1532,1).Value = 1531
$ws.Cells.Item(1532,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1532,3).Value = "4:05 PM"
$ws.Cells.Item(1532,4).Value = "W61475"
$ws.Cells.Item(1532,5).Value = "Barcelona"
$ws.Cells.Item(1532,6).Value = "(BCN)"
$ws.Cells.Item(1532,7).Value = "Wizz Air "
$ws.Cells.Item(1532,8).Value = "A321"
$ws.Cells.Item(1532,9).Value = "(HA-LXG)"
$ws.Cells.Item(1532,10).Value = "5:18 PM"
$ws.Cells.Item(1532,12).Value = "1 hours, 13 minutes"

# Row 1533
$ws.Cells.Item(1533,1).Value = 1532
$ws.Cells.Item(1533,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1533,3).Value = "4:10 PM"
$ws.Cells.Item(1533,4).Value = "LO71"
$ws.Cells.Item(1533,5).Value = "Delhi"
$ws.Cells.Item(1533,6).Value = "(DEL)"
$ws.Cells.Item(1533,7).Value = "LOT (Independence Livery) "
$ws.Cells.Item(1533,8).Value = "B789"
$ws.Cells.Item(1533,9).Value = "(SP-LSC)"
$ws.Cells.Item(1533,10).Value = "4:31 PM"
$ws.Cells.Item(1533,12).Value = "0 hours, 21 minutes"

# Row 1534
$ws.Cells.Item(1534,1).Value = 1533
$ws.Cells.Item(1534,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1534,3).Value = "4:15 PM"
$ws.Cells.Item(1534,4).Value = "LO343"
$ws.Cells.Item(1534,5).Value = "Nice"
$ws.Cells.Item(1534,6).Value = "(NCE)"
$ws.Cells.Item(1534,7).Value = "LOT (Star Alliance Livery) "
$ws.Cells.Item(1534,8).Value = "E75S"
$ws.Cells.Item(1534,9).Value = "(SP-LIO)"
$ws.Cells.Item(1534,10).Value = "5:06 PM"
$ws.Cells.Item(1534,12).Value = "0 hours, 51 minutes"

# Row 1535
$ws.Cells.Item(1535,1).Value = 1534
$ws.Cells.Item(1535,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1535,3).Value = "4:25 PM"
$ws.Cells.Item(1535,4).Value = "LO335"
$ws.Cells.Item(1535,5).Value = "Paris"
$ws.Cells.Item(1535,6).Value = "(CDG)"
$ws.Cells.Item(1535,7).Value = "LOT "
$ws.Cells.Item(1535,8).Value = "E195"
$ws.Cells.Item(1535,9).Value = "(SP-LNI)"
$ws.Cells.Item(1535,10).Value = "4:44 PM"
$ws.Cells.Item(1535,12).Value = "0 hours, 19 minutes"

# Row 1536
$ws.Cells.Item(1536,1).Value = 1535
$ws.Cells.Item(1536,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1536,3).Value = "4:25 PM"
$ws.Cells.Item(1536,4).Value = "LO775"
$ws.Cells.Item(1536,5).Value = "Vilnius"
$ws.Cells.Item(1536,6).Value = "(VNO)"
$ws.Cells.Item(1536,7).Value = "LOT "
$ws.Cells.Item(1536,8).Value = "E195"
$ws.Cells.Item(1536,9).Value = "(SP-LNG)"
$ws.Cells.Item(1536,10).Value = "4:38 PM"
$ws.Cells.Item(1536,12).Value = "0 hours, 13 minutes"

# Row 1537
$ws.Cells.Item(1537,1).Value = 1536
$ws.Cells.Item(1537,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1537,3).Value = "4:35 PM"
$ws.Cells.Item(1537,4).Value = "LO233"
$ws.Cells.Item(1537,5).Value = "Brussels"
$ws.Cells.Item(1537,6).Value = "(BRU)"
$ws.Cells.Item(1537,7).Value = "LOT "
$ws.Cells.Item(1537,8).Value = "E195"
$ws.Cells.Item(1537,9).Value = "(SP-LNM)"
$ws.Cells.Item(1537,10).Value = "4:53 PM"
$ws.Cells.Item(1537,12).Value = "0 hours, 18 minutes"

# Row 1538
$ws.Cells.Item(1538,1).Value = 1537
$ws.Cells.Item(1538,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1538,3).Value = "4:35 PM"
$ws.Cells.Item(1538,4).Value = "LO389"
$ws.Cells.Item(1538,5).Value = "Berlin"
$ws.Cells.Item(1538,6).Value = "(BER)"
$ws.Cells.Item(1538,7).Value = "LOT "
$ws.Cells.Item(1538,8).Value = "E170"
$ws.Cells.Item(1538,9).Value = "(SP-LDI)"
$ws.Cells.Item(1538,10).Value = "4:51 PM"
$ws.Cells.Item(1538,12).Value = "0 hours, 16 minutes"

# Row 1539
$ws.Cells.Item(1539,1).Value = 1538
$ws.Cells.Item(1539,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1539,3).Value = "4:35 PM"
$ws.Cells.Item(1539,4).Value = "LO3923"
$ws.Cells.Item(1539,5).Value = "Krakow"
$ws.Cells.Item(1539,6).Value = "(KRK)"
$ws.Cells.Item(1539,7).Value = "LOT "
$ws.Cells.Item(1539,8).Value = "E190"
$ws.Cells.Item(1539,9).Value = "(SP-LME)"
$ws.Cells.Item(1539,10).Value = "5:00 PM"
$ws.Cells.Item(1539,12).Value = "0 hours, 25 minutes"

# Row 1540
$ws.Cells.Item(1540,1).Value = 1539
$ws.Cells.Item(1540,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1540,3).Value = "4:40 PM"
$ws.Cells.Item(1540,4).Value = "LO3"
$ws.Cells.Item(1540,5).Value = "Chicago"
$ws.Cells.Item(1540,6).Value = "(ORD)"
$ws.Cells.Item(1540,7).Value = "LOT "
$ws.Cells.Item(1540,8).Value = "B788"
$ws.Cells.Item(1540,9).Value = "(SP-LRF)"
$ws.Cells.Item(1540,10).Value = "5:12 PM"
$ws.Cells.Item(1540,12).Value = "0 hours, 32 minutes"

# Row 1541
$ws.Cells.Item(1541,1).Value = 1540
$ws.Cells.Item(1541,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1541,3).Value = "4:40 PM"
$ws.Cells.Item(1541,4).Value = "LO267"
$ws.Cells.Item(1541,5).Value = "Amsterdam"
$ws.Cells.Item(1541,6).Value = "(AMS)"
$ws.Cells.Item(1541,7).Value = "LOT "
$ws.Cells.Item(1541,8).Value = "E195"
$ws.Cells.Item(1541,9).Value = "(SP-LNO)"
$ws.Cells.Item(1541,10).Value = "4:58 PM"
$ws.Cells.Item(1541,12).Value = "0 hours, 18 minutes"

# Row 1542
$ws.Cells.Item(1542,1).Value = 1541
$ws.Cells.Item(1542,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1542,3).Value = "4:40 PM"
$ws.Cells.Item(1542,4).Value = "LO415"
$ws.Cells.Item(1542,5).Value = "Geneva"
$ws.Cells.Item(1542,6).Value = "(GVA)"
$ws.Cells.Item(1542,7).Value = "LOT "
$ws.Cells.Item(1542,8).Value = "B738"
$ws.Cells.Item(1542,9).Value = "(SP-LWD)"
$ws.Cells.Item(1542,10).Value = "4:57 PM"
$ws.Cells.Item(1542,12).Value = "0 hours, 17 minutes"

# Row 1543
$ws.Cells.Item(1543,1).Value = 1542
$ws.Cells.Item(1543,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1543,3).Value = "4:50 PM"
$ws.Cells.Item(1543,4).Value = "LO15"
$ws.Cells.Item(1543,5).Value = "New York"
$ws.Cells.Item(1543,6).Value = "(EWR)"
$ws.Cells.Item(1543,7).Value = "LOT "
$ws.Cells.Item(1543,8).Value = "B788"
$ws.Cells.Item(1543,9).Value = "(SP-LRH)"
$ws.Cells.Item(1543,10).Value = "5:40 PM"
$ws.Cells.Item(1543,12).Value = "0 hours, 50 minutes"

# Row 1544
$ws.Cells.Item(1544,1).Value = 1543
$ws.Cells.Item(1544,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1544,3).Value = "4:50 PM"
$ws.Cells.Item(1544,4).Value = "LO26"
$ws.Cells.Item(1544,5).Value = "New York"
$ws.Cells.Item(1544,6).Value = "(JFK)"
$ws.Cells.Item(1544,7).Value = "LOT "
$ws.Cells.Item(1544,8).Value = "B788"
$ws.Cells.Item(1544,9).Value = "(SP-LRB)"
$ws.Cells.Item(1544,10).Value = "5:43 PM"
$ws.Cells.Item(1544,12).Value = "0 hours, 53 minutes"

# Row 1545
$ws.Cells.Item(1545,1).Value = 1544
$ws.Cells.Item(1545,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1545,3).Value = "4:50 PM"
$ws.Cells.Item(1545,4).Value = "LO319"
$ws.Cells.Item(1545,5).Value = "Milan"
$ws.Cells.Item(1545,6).Value = "(MXP)"
$ws.Cells.Item(1545,7).Value = "LOT "
$ws.Cells.Item(1545,8).Value = "E190"
$ws.Cells.Item(1545,9).Value = "(SP-LMH)"
$ws.Cells.Item(1545,10).Value = "5:04 PM"
$ws.Cells.Item(1545,12).Value = "0 hours, 14 minutes"

# Row 1546
$ws.Cells.Item(1546,1).Value = 1545
$ws.Cells.Item(1546,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1546,3).Value = "4:50 PM"
$ws.Cells.Item(1546,4).Value = "LO483"
$ws.Cells.Item(1546,5).Value = "Oslo"
$ws.Cells.Item(1546,6).Value = "(OSL)"
$ws.Cells.Item(1546,7).Value = "LOT "
$ws.Cells.Item(1546,8).Value = "B738"
$ws.Cells.Item(1546,9).Value = "(SP-LWC)"
$ws.Cells.Item(1546,10).Value = "5:37 PM"
$ws.Cells.Item(1546,12).Value = "0 hours, 47 minutes"

# Row 1547
$ws.Cells.Item(1547,1).Value = 1546
$ws.Cells.Item(1547,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1547,3).Value = "4:50 PM"
$ws.Cells.Item(1547,4).Value = "LO3947"
$ws.Cells.Item(1547,5).Value = "Poznan"
$ws.Cells.Item(1547,6).Value = "(POZ)"
$ws.Cells.Item(1547,7).Value = "LOT "
$ws.Cells.Item(1547,8).Value = "E75S"
$ws.Cells.Item(1547,9).Value = "(SP-LID)"
$ws.Cells.Item(1547,10).Value = "5:03 PM"
$ws.Cells.Item(1547,12).Value = "0 hours, 13 minutes"

# Row 1548
$ws.Cells.Item(1548,1).Value = 1547
$ws.Cells.Item(1548,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1548,3).Value = "5:00 PM"
$ws.Cells.Item(1548,4).Value = "LH1615"
$ws.Cells.Item(1548,5).Value = "Munich"
$ws.Cells.Item(1548,6).Value = "(MUC)"
$ws.Cells.Item(1548,7).Value = "Lufthansa "
$ws.Cells.Item(1548,8).Value = "CRJ9"
$ws.Cells.Item(1548,9).Value = "(D-ACNX)"
$ws.Cells.Item(1548,10).Value = "5:01 PM"
$ws.Cells.Item(1548,12).Value = "0 hours, 1 minutes"

# Row 1549
$ws.Cells.Item(1549,1).Value = 1548
$ws.Cells.Item(1549,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1549,3).Value = "5:00 PM"
$ws.Cells.Item(1549,4).Value = "LO45"
$ws.Cells.Item(1549,5).Value = "Toronto"
$ws.Cells.Item(1549,6).Value = "(YYZ)"
$ws.Cells.Item(1549,7).Value = "LOT "
$ws.Cells.Item(1549,8).Value = "B789"
$ws.Cells.Item(1549,9).Value = "(SP-LSB)"
$ws.Cells.Item(1549,10).Value = "5:34 PM"
$ws.Cells.Item(1549,12).Value = "0 hours, 34 minutes"

# Row 1550
$ws.Cells.Item(1550,1).Value = 1549
$ws.Cells.Item(1550,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1550,3).Value = "5:00 PM"
$ws.Cells.Item(1550,4).Value = "LO353"
$ws.Cells.Item(1550,5).Value = "Munich"
$ws.Cells.Item(1550,6).Value = "(MUC)"
$ws.Cells.Item(1550,7).Value = "LOT "
$ws.Cells.Item(1550,8).Value = "E190"
$ws.Cells.Item(1550,9).Value = "(SP-LMA)"
$ws.Cells.Item(1550,10).Value = "5:16 PM"
$ws.Cells.Item(1550,12).Value = "0 hours, 16 minutes"

# Row 1551
$ws.Cells.Item(1551,1).Value = 1550
$ws.Cells.Item(1551,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1551,3).Value = "5:00 PM"
$ws.Cells.Item(1551,4).Value = "LO379"
$ws.Cells.Item(1551,5).Value = "Frankfurt"
$ws.Cells.Item(1551,6).Value = "(FRA)"
$ws.Cells.Item(1551,7).Value = "LOT "
$ws.Cells.Item(1551,8).Value = "E190"
$ws.Cells.Item(1551,9).Value = "(SP-LMB)"
$ws.Cells.Item(1551,10).Value = "5:14 PM"
$ws.Cells.Item(1551,12).Value = "0 hours, 14 minutes"

# Row 1552
$ws.Cells.Item(1552,1).Value = 1551
$ws.Cells.Item(1552,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1552,3).Value = "5:00 PM"
$ws.Cells.Item(1552,4).Value = "LO419"
$ws.Cells.Item(1552,5).Value = "Zurich"
$ws.Cells.Item(1552,6).Value = "(ZRH)"
$ws.Cells.Item(1552,7).Value = "LOT (Independence Livery) "
$ws.Cells.Item(1552,8).Value = "B38M"
$ws.Cells.Item(1552,9).Value = "(SP-LVD)"
$ws.Cells.Item(1552,10).Value = "5:21 PM"
$ws.Cells.Item(1552,12).Value = "0 hours, 21 minutes"

# Row 1553
$ws.Cells.Item(1553,1).Value = 1552
$ws.Cells.Item(1553,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1553,3).Value = "5:00 PM"
$ws.Cells.Item(1553,4).Value = "LO457"
$ws.Cells.Item(1553,5).Value = "Stockholm"
$ws.Cells.Item(1553,6).Value = "(ARN)"
$ws.Cells.Item(1553,7).Value = "LOT "
$ws.Cells.Item(1553,8).Value = "B38M"
$ws.Cells.Item(1553,9).Value = "(SP-LVB)"
$ws.Cells.Item(1553,10).Value = "5:45 PM"
$ws.Cells.Item(1553,12).Value = "0 hours, 45 minutes"

# Row 1554
$ws.Cells.Item(1554,1).Value = 1553
$ws.Cells.Item(1554,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1554,3).Value = "5:10 PM"
$ws.Cells.Item(1554,4).Value = "LO373"
$ws.Cells.Item(1554,5).Value = "Stuttgart"
$ws.Cells.Item(1554,6).Value = "(STR)"
$ws.Cells.Item(1554,7).Value = "LOT "
$ws.Cells.Item(1554,8).Value = "E195"
$ws.Cells.Item(1554,9).Value = "(SP-LNK)"
$ws.Cells.Item(1554,10).Value = "5:20 PM"
$ws.Cells.Item(1554,12).Value = "0 hours, 10 minutes"

# Row 1555
$ws.Cells.Item(1555,1).Value = 1554
$ws.Cells.Item(1555,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1555,3).Value = "5:10 PM"
$ws.Cells.Item(1555,4).Value = "LO407"
$ws.Cells.Item(1555,5).Value = "Dusseldorf"
$ws.Cells.Item(1555,6).Value = "(DUS)"
$ws.Cells.Item(1555,7).Value = "LOT "
$ws.Cells.Item(1555,8).Value = "E75S"
$ws.Cells.Item(1555,9).Value = "(SP-LIK)"
$ws.Cells.Item(1555,10).Value = "5:29 PM"
$ws.Cells.Item(1555,12).Value = "0 hours, 19 minutes"

# Row 1556
$ws.Cells.Item(1556,1).Value = 1555
$ws.Cells.Item(1556,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1556,3).Value = "5:10 PM"
$ws.Cells.Item(1556,4).Value = "LO393"
$ws.Cells.Item(1556,5).Value = "Hamburg"
$ws.Cells.Item(1556,6).Value = "(HAM)"
$ws.Cells.Item(1556,7).Value = "LOT "
$ws.Cells.Item(1556,8).Value = "E75S"
$ws.Cells.Item(1556,9).Value = "(SP-LIB)"
$ws.Cells.Item(1556,10).Value = "5:25 PM"
$ws.Cells.Item(1556,12).Value = "0 hours, 15 minutes"

# Row 1557
$ws.Cells.Item(1557,1).Value = 1556
$ws.Cells.Item(1557,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1557,3).Value = "5:20 PM"
$ws.Cells.Item(1557,4).Value = "AF1047"
$ws.Cells.Item(1557,5).Value = "Paris"
$ws.Cells.Item(1557,6).Value = "(CDG)"
$ws.Cells.Item(1557,7).Value = "Air France "
$ws.Cells.Item(1557,8).Value = "E190"
$ws.Cells.Item(1557,9).Value = "(F-HBLK)"
$ws.Cells.Item(1557,10).Value = "5:35 PM"
$ws.Cells.Item(1557,12).Value = "0 hours, 15 minutes"

# Row 1558
$ws.Cells.Item(1558,1).Value = 1557
$ws.Cells.Item(1558,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1558,3).Value = "5:20 PM"
$ws.Cells.Item(1558,4).Value = "LO149"
$ws.Cells.Item(1558,5).Value = "Cairo"
$ws.Cells.Item(1558,6).Value = "(CAI)"
$ws.Cells.Item(1558,7).Value = "LOT "
$ws.Cells.Item(1558,8).Value = "E190"
$ws.Cells.Item(1558,9).Value = "(SP-LMF)"
$ws.Cells.Item(1558,10).Value = "5:46 PM"
$ws.Cells.Item(1558,12).Value = "0 hours, 26 minutes"

# Row 1559
$ws.Cells.Item(1559,1).Value = 1558
$ws.Cells.Item(1559,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1559,3).Value = "5:20 PM"
$ws.Cells.Item(1559,4).Value = "LO253"
$ws.Cells.Item(1559,5).Value = "Luxembourg"
$ws.Cells.Item(1559,6).Value = "(LUX)"
$ws.Cells.Item(1559,7).Value = "LOT (Sliwka Naleczowska Livery) "
$ws.Cells.Item(1559,8).Value = "E195"
$ws.Cells.Item(1559,9).Value = "(SP-LNC)"
$ws.Cells.Item(1559,10).Value = "5:32 PM"
$ws.Cells.Item(1559,12).Value = "0 hours, 12 minutes"

# Row 1560
$ws.Cells.Item(1560,1).Value = 1559
$ws.Cells.Item(1560,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1560,3).Value = "5:25 PM"
$ws.Cells.Item(1560,4).Value = "LO461"
$ws.Cells.Item(1560,5).Value = "Copenhagen"
$ws.Cells.Item(1560,6).Value = "(CPH)"
$ws.Cells.Item(1560,7).Value = "LOT "
$ws.Cells.Item(1560,8).Value = "E195"
$ws.Cells.Item(1560,9).Value = "(SP-LNP)"
$ws.Cells.Item(1560,10).Value = "5:49 PM"
$ws.Cells.Item(1560,12).Value = "0 hours, 24 minutes"

# Row 1561
$ws.Cells.Item(1561,1).Value = 1560
$ws.Cells.Item(1561,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1561,3).Value = "5:30 PM"
$ws.Cells.Item(1561,4).Value = "LO467"
$ws.Cells.Item(1561,5).Value = "Billund"
$ws.Cells.Item(1561,6).Value = "(BLL)"
$ws.Cells.Item(1561,7).Value = "LOT "
$ws.Cells.Item(1561,8).Value = "E75S"
$ws.Cells.Item(1561,9).Value = "(SP-LIN)"
$ws.Cells.Item(1561,10).Value = "5:51 PM"
$ws.Cells.Item(1561,12).Value = "0 hours, 21 minutes"

# Row 1562
$ws.Cells.Item(1562,1).Value = 1561
$ws.Cells.Item(1562,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1562,3).Value = "5:30 PM"
$ws.Cells.Item(1562,4).Value = "LPR42"
$ws.Cells.Item(1562,5).Value = "Gdansk"
$ws.Cells.Item(1562,6).Value = "(GDN)"
$ws.Cells.Item(1562,7).Value = "Polish Medical Air Rescue "
$ws.Cells.Item(1562,8).Value = "LJ75"
$ws.Cells.Item(1562,9).Value = "(SP-MXS)"
$ws.Cells.Item(1562,10).Value = "5:48 PM"
$ws.Cells.Item(1562,12).Value = "0 hours, 18 minutes"

# Row 1563
$ws.Cells.Item(1563,1).Value = 1562
$ws.Cells.Item(1563,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1563,3).Value = "5:40 PM"
$ws.Cells.Item(1563,4).Value = "A3873"
$ws.Cells.Item(1563,5).Value = "Athens"
$ws.Cells.Item(1563,6).Value = "(ATH)"
$ws.Cells.Item(1563,7).Value = "Aegean Airlines "
$ws.Cells.Item(1563,8).Value = "A320"
$ws.Cells.Item(1563,9).Value = "(SX-DVR)"
$ws.Cells.Item(1563,10).Value = "6:06 PM"
$ws.Cells.Item(1563,12).Value = "0 hours, 26 minutes"

# Row 1564
$ws.Cells.Item(1564,1).Value = 1563
$ws.Cells.Item(1564,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(1564,3).Value = "5:50 PM"
$ws.Cells.Item(1564,4).Value = "LO525"
$ws.Cells.Item(1564,5).Value = "Prague"
$ws.Cells.Item(1564,6).Value = "(PRG)"
$ws.Cells.Item(1564,7).Value = "LOT "
$ws.Cells.Item(1564,8).Value = "E195"
$ws.Cells.Item(1564,9).Value = "(SP-LNH)"
$ws.Cells.Item(1564,10).Value = "6:02 PM"
$ws.Cells.Item(1564,12).Value = "0 hours, 12 minutes"
